$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new blank row above row 2 (shifts existing rows down)
$ws.Rows.Item(2).Insert()

# Merge the new row A2:H2 and apply the same styling as the title row (A1:H1)
$ws.Range("A1:H1").Copy()
$ws.Range("A2:H2").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A2:H2").Merge()
$ws.Rows.Item(2).RowHeight = 22.05

$ws.Range("K3").Select()
